$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Update browser values for the amazonHamburgerMenuTest rows (rows 8 and 9)
# so they use "chrome" instead of "edge"/"firefox".
$ws.Range("C8").Value = "chrome"
$ws.Range("C9").Value = "chrome"

# Move the active selection from B8 to C7.
$ws.Activate()
$ws.Range("C7").Select()
